# Update "想去人数" (number of people interested) counts that changed
# when the source data was re-scraped.
#
#   Sheet 展览   (Exhibitions):      F6 5161->5162, F16 4229->4230, F29 1078->1079, F31 2707->2709
#   Sheet 全部类型 (All types):       F6 5161->5162, F16 4229->4230, F30 1078->1079, F32 2707->2709

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 5162
$wsExhibit.Range("F16").Value = 4230
$wsExhibit.Range("F29").Value = 1079
$wsExhibit.Range("F31").Value = 2709

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 5162
$wsAll.Range("F16").Value = 4230
$wsAll.Range("F30").Value = 1079
$wsAll.Range("F32").Value = 2709
